# Applies the commit's edits to analise_integracao.docx:
#  1. Clear the "Processos desejados:" and "Informações necessárias pelo ERP:"
#     values on the second paragraph (keep the labels, drop the content).
#  2. Replace the long "Analise Funcional Recomendada" body (scenario table +
#     JSON sample + closing remarks) with four short summary sentences.
#  3. Strip every data row from the "Mapeamento de Campos ERP" table, keeping
#     only the header row (Campo | Descricao).

$d = $word.ActiveDocument

# --- 1. Second paragraph: drop the values after the two labels ---
$d.Content.Find.Execute('Processos desejados: Reembolso', $true, $false, $false, $false, $false, $true, 1, $false, 'Processos desejados: ', 2) | Out-Null
$d.Content.Find.Execute('Informações necessárias pelo ERP: Valor total relatório, CPF, Empresa, Filial, Divisão, Motivo, Tipo de documento, Rateio, Valor do rateio, Centro de custo, PEP', $true, $false, $false, $false, $false, $true, 1, $false, 'Informações necessárias pelo ERP: ', 2) | Out-Null

# --- 2. Fourth paragraph: replace the whole narrative block (spans several
#        manual line breaks -- Chr(11) below stands for <w:br/>) with the
#        four new summary sentences. ---
$d.Content.Find.Execute('Para a integração do ERP SAP ECC/4HANA com o SaaS Paytrack, a análise funcional pode ser realizada da seguinte forma:Cenário: ReembolsoMapeamento de Campos:| Campo ERP SAP ECC/4HANA | Campo Paytrack ||------------------------|----------------|| BUKRS (Empresa)        | Empresa        || WERKS (Filial)         | Filial         || SPART (Divisão)        | Divisão        || ZFBDT (Motivo)         | Motivo         || BLART (Tipo de documento) | Tipo de documento || KOSTL (Centro de custo) | Centro de custo || PRCTR (PEP)            | PEP            || LFBNR (Valor total do relatório) | Valor total do relatório || STCD1 (CPF)            | CPF            || ZPRCTR (Rateio)        | Rateio         || ZKOSTL (Valor do rateio) | Valor do rateio |Exemplo de JSON de Retorno:```{  "BUKRS": "1234",  "WERKS": "5678",  "SPART": "9ABC",  "ZFBDT": "Motivo do reembolso",  "BLART": "Tipo de documento",  "KOSTL": "Centro de custo",  "PRCTR": "PEP",  "LFBNR": "1000.00",  "STCD1": "123.456.789-00",  "ZPRCTR": "Rateio",  "ZKOSTL": "500.00"}```É importante destacar que a análise funcional para outros cenários, como adiantamento e prestação de contas, deve ser realizada de forma separada, seguindo o mesmo padrão.Além disso, é necessário que o cliente disponibilize um Webservice para que a Paytrack possa consumir as informações necessárias para a integração. A comunicação será realizada de forma sincrona, com a Paytrack ativa nas integrações.', $true, $false, $false, $false, $false, $true, 1, $false, 'Uma boa abordagem para realizar a integração do ERP SAP ECC/4HANA com o SaaS Paytrack seria realizar uma análise funcional detalhada dos processos desejados e das informações necessárias pelo ERP repassadas pelo cliente. Nesta análise funcional, é importante incluir um mapeamento de campos olhando para o ERP selecionado, apresentando as nomenclaturas específicas do ERP, como bukrs para empresa, em um formato de tabela para facilitar a compreensão. Além disso, é essencial incluir um exemplo de JSON formatado com as nomenclaturas do ERP para ilustrar como as informações devem ser estruturadas.Para garantir uma integração eficiente, é fundamental seguir algumas diretrizes, como utilizar comunicação síncrona com os webservices do cliente, garantir que a Paytrack seja ativa nas integrações e separar a análise funcional por cenário selecionado, como adiantamento e prestação de contas, para facilitar a implementação e garantir que todos os requisitos sejam atendidos de forma adequada. Dessa forma, uma análise funcional detalhada e organizada será essencial para garantir o sucesso da integração entre o ERP SAP ECC/4HANA e o SaaS Paytrack.', 2) | Out-Null

# --- 3. Table: drop every data row, keep only the header row ---
$tbl = $d.Tables.Item(1)
for ($i = $tbl.Rows.Count; $i -ge 2; $i--) {
    $tbl.Rows.Item($i).Delete()
}

Write-Host "Edits applied."
